$wb = $excel.ActiveWorkbook

# ---- Big prompt text blocks (here-strings) ----
$s5 = @"
 Given is the adjacency matrix for a weighted undirected graph containing 16 nodes labelled A to P. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   

Consider some examples

Example 1: what is the least cost path from node A to node I? Return the sequence of nodes in response.
   A B C D E F G H I
 A 0 1 0 3 0 0 0 0 0
 B 1 0 2 0 2 0 0 0 0
 C 0 2 0 0 0 2 0 0 0
 D 3 0 0 0 1 0 2 0 0
 E 0 2 0 1 0 3 0 1 0
 F 0 0 2 0 3 0 0 0 1
 G 0 0 0 2 0 0 0 2 0
 H 0 0 0 0 1 0 2 0 1
 I 0 0 0 0 0 1 0 1 0

Solution: A -> B -> E -> H -> I
        

Example 2: what is the least cost path from node A to node I? Return the sequence of nodes in response.
   A B C D E F G H I
 A 0 4 0 2 0 0 0 0 0
 B 4 0 4 0 3 0 0 0 0
 C 0 4 0 0 0 3 0 0 0
 D 2 0 0 0 0 0 4 0 0
 E 0 3 0 0 0 3 0 0 0
 F 0 0 3 0 3 0 0 0 1
 G 0 0 0 4 0 0 0 3 0
 H 0 0 0 0 0 0 3 0 4
 I 0 0 0 0 0 1 0 4 0

Solution: A -> B -> E -> F -> I
        

Example 3: what is the least cost path from node A to node P? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P
 A 0 4 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 B 4 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 2 0 5 0 0 4 0 0 0 0 0 0 0 0 0
 D 0 0 5 0 0 0 0 3 0 0 0 0 0 0 0 0
 E 1 0 0 0 0 2 0 0 4 0 0 0 0 0 0 0
 F 0 0 0 0 2 0 4 0 0 0 0 0 0 0 0 0
 G 0 0 4 0 0 4 0 5 0 0 5 0 0 0 0 0
 H 0 0 0 3 0 0 5 0 0 0 0 5 0 0 0 0
 I 0 0 0 0 4 0 0 0 0 2 0 0 4 0 0 0
 J 0 0 0 0 0 0 0 0 2 0 0 0 0 4 0 0
 K 0 0 0 0 0 0 5 0 0 0 0 2 0 0 1 0
 L 0 0 0 0 0 0 0 5 0 0 2 0 0 0 0 4
 M 0 0 0 0 0 0 0 0 4 0 0 0 0 5 0 0
 N 0 0 0 0 0 0 0 0 0 4 0 0 5 0 5 0
 O 0 0 0 0 0 0 0 0 0 0 1 0 0 5 0 3
 P 0 0 0 0 0 0 0 0 0 0 0 4 0 0 3 0

Solution: A -> E -> F -> G -> K -> O -> P
        
 Given these examples, answer the following quesiton.

what is the least cost path from node A to node P? Return the sequence of nodes in response.

   A B C D E F G H I J K L M N O P
 A 0 5 0 0 5 0 0 0 0 0 0 0 0 0 0 0
 B 5 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 4 0 5 0 0 5 0 0 0 0 0 0 0 0 0
 D 0 0 5 0 0 0 0 5 0 0 0 0 0 0 0 0
 E 5 0 0 0 0 5 0 0 3 0 0 0 0 0 0 0
 F 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 5 0 0 0 0 1 0 0 0 0 0 0 0 0
 H 0 0 0 5 0 0 1 0 0 0 0 3 0 0 0 0
 I 0 0 0 0 3 0 0 0 0 4 0 0 2 0 0 0
 J 0 0 0 0 0 0 0 0 4 0 0 0 0 5 0 0
 K 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0
 L 0 0 0 0 0 0 0 3 0 0 4 0 0 0 0 1
 M 0 0 0 0 0 0 0 0 2 0 0 0 0 4 0 0
 N 0 0 0 0 0 0 0 0 0 5 0 0 4 0 5 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 5
 P 0 0 0 0 0 0 0 0 0 0 0 1 0 0 5 0
    
"@
$s10 = @"
 Given is the adjacency matrix for a weighted undirected graph containing 23 nodes labelled A to W. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   
Consider some examples
Example 1: what is the least cost path from node A to node Y? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 4 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 4 0 2 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 2 0 1 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 4 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 4 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 2 0 0 0 0 0 2 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 1 0 0 0 2 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 3 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 1 0 0 0 4 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 4 0 0 0 3 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 5 0 0 0 0 0 5 0 0 0 2 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 5 0 2 0 0 0 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 4 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 3 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 5 0 0 0 1 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 1 0 0 0 3 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 1 0 1 0 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 2
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 2 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 2 0 2 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 2 0 2 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 2
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 2 0
Solution: A -> B -> C -> D -> I -> J -> O -> T -> Y
Example 2: what is the least cost path from node A to node Y? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 2 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 2 0 2 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 2 0 3 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 3 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 5 0 0 0 0 0 3 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 1 0 0 0 3 0 5 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 2 0 0 0 5 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 2 0 0 0 1 0 3 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 1 0 0 0 3 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 3 0 0 0 0 0 4 0 0 0 3 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 3 0 0 0 4 0 0 0 0 0 3 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 3 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 3 0 0 0 4 0 0 0 0 0 5 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 1 0 0 0 3 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 3 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 3
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 4 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 4 0 5 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 5 0 1 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 1 0 5
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 5 0
Solution: A -> B -> G -> L -> Q -> R -> W -> X -> Y
Example 3: what is the least cost path from node A to node X? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X
 A 0 5 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 5 0 2 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 2 0 3 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 3 0 1 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 2 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 1 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 4 0 0 0 3 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 4 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 5 0 0 0 1 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 1 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 5 0 0 0 1 0 0 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 0 2 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 1 0 0 0 5 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 2 0 0 0 5 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 4 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 5
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 3 0 0 0
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 4 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 4 0 5 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 5 0 4
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 4 0
Solution: A -> F -> K -> O -> T -> U -> V -> W -> X
 Given these examples, answer the following quesiton.
what is the least cost path from node A to node W? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W
 A 0 2 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 2 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 3 0 5 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 5 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 1 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 3 0 0 0 0 1 0 0 4 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 1 0 0 1 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 3 0 0 0 0 1 0 0 5 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 3 0 0 1 0 0 0 0 2 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0 0 5 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 5 0 0 0 0 3 0 0 0 1 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 2 0 0 3 0 3 0 0 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 1 0 0 0 3 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 1 0 5 0 0 0 3 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 5 0 0 0 0 0 5
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 4 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 5 0 0
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 5 0 5 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 5 0 5
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 5 0
    
"@
$s14 = @"
 Given is the adjacency matrix for a weighted undirected graph containing 24 nodes labelled A to X. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   
Consider some examples
Example 1: what is the least cost path from node A to node Y? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 4 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 4 0 2 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 2 0 1 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 4 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 4 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 2 0 0 0 0 0 2 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 1 0 0 0 2 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 3 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 1 0 0 0 4 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 4 0 0 0 3 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 5 0 0 0 0 0 5 0 0 0 2 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 5 0 2 0 0 0 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 4 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 3 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 5 0 0 0 1 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 1 0 0 0 3 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 1 0 1 0 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 2
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 2 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 2 0 2 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 2 0 2 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 2
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 2 0
Solution: A -> B -> C -> D -> I -> J -> O -> T -> Y
Example 2: what is the least cost path from node A to node Y? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 2 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 2 0 2 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 2 0 3 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 3 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 5 0 0 0 0 0 3 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 1 0 0 0 3 0 5 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 2 0 0 0 5 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 2 0 0 0 1 0 3 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 1 0 0 0 3 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 3 0 0 0 0 0 4 0 0 0 3 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 3 0 0 0 4 0 0 0 0 0 3 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 3 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 3 0 0 0 4 0 0 0 0 0 5 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 1 0 0 0 3 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 3 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 3
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 4 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 4 0 5 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 5 0 1 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 1 0 5
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 5 0
Solution: A -> B -> G -> L -> Q -> R -> W -> X -> Y
Example 3: what is the least cost path from node A to node X? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X
 A 0 5 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 5 0 2 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 2 0 3 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 3 0 1 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 2 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 1 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 4 0 0 0 3 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 4 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 5 0 0 0 1 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 1 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 5 0 0 0 1 0 0 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 0 2 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 1 0 0 0 5 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 2 0 0 0 5 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 4 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 5
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 3 0 0 0
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 4 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 4 0 5 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 5 0 4
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 4 0
Solution: A -> F -> K -> O -> T -> U -> V -> W -> X
 Given these examples, answer the following quesiton.
what is the least cost path from node A to node X? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X
 A 0 2 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 2 0 5 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 5 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 4 0 4 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 4 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 3 0 0 0 0 0 3 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 1 0 0 0 3 0 4 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 4 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 3 0 0 0 3 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 3 0 0 0 1 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 1 0 0 0 0 0 1 0 0 4 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 3 0 0 0 1 0 0 0 0 2 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 0 4 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 3 0 0 0 1 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 2 0 0 3 0 4 0 0 0 1 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 1 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 1 0 2 0 0 0 3 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 2 0 0 0 0 0 5
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 3 0 0 0
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 3 0 4 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 3 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 3 0 5
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 5 0
    
"@

# ---- Sheet 1 (o_10): add column E (evaluator_partial_correctness) and update row 2 ----
$ws1 = $wb.Worksheets.Item("o_10")
$ws1.Range("E1").Value = "evaluator_partial_correctness"
$ws1.Range("A2").Value = $s5
$ws1.Range("B2").Value = "A -> B -> C -> G -> H -> L -> P"
$ws1.Range("C2").Value = "Solution: A -> E -> I -> M -> N -> O -> P"
$ws1.Range("D2").Value = "invalid input"
$ws1.Range("E2").Value = "1/7"

# ---- Sheet 2 (o_20): new sheet, added after o_10 ----
$after1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $after1)
$ws2.Name = "o_20"
$ws2.Range("A1").Value = "prompt"
$ws2.Range("B1").Value = "solution"
$ws2.Range("C1").Value = "llm_response"
$ws2.Range("D1").Value = "evaluator_response"
$ws2.Range("E1").Value = "evaluator_partial_correctness"
$ws2.Range("A2").Value = $s10
$ws2.Range("B2").Value = "A -> F -> J -> K -> O -> P -> Q -> V -> W"
$ws2.Range("C2").Value = "Solution: A -> F -> J -> K -> O -> P -> Q -> V -> W"
$ws2.Range("D2").Value = "invalid input"
$ws2.Range("E2").Value = "9/9"

# ---- Sheet 3 (o_20_jumbled): new sheet, added after o_20 ----
$after2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Type]::Missing, $after2)
$ws3.Name = "o_20_jumbled"
$ws3.Range("A1").Value = "prompt"
$ws3.Range("B1").Value = "solution"
$ws3.Range("C1").Value = "llm_response"
$ws3.Range("D1").Value = "evaluator_response"
$ws3.Range("E1").Value = "evaluator_partial_correctness"
$ws3.Range("A2").Value = $s14
$ws3.Range("B2").Value = "A -> F -> K -> L -> P -> Q -> R -> S -> X"
$ws3.Range("C2").Value = "Solution: A -> B -> G -> L -> P -> U -> V -> W -> X"
$ws3.Range("D2").Value = "invalid input"
$ws3.Range("E2").Value = "0/9"

# ---- Restore o_10 as the active/selected sheet ----
$ws1.Activate()

"Worksheets: " + $wb.Worksheets.Count
for ($i=1; $i -le $wb.Worksheets.Count; $i++) { "Sheet " + $i + ": " + $wb.Worksheets.Item($i).Name }
